$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2020" column (Q) added to the growth-rate table, mirroring the
# formatting already used by the preceding "2019" column (P).

# Header row: year label
$ws.Range("Q4").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)  # xlPasteFormats

# Data row: growth-rate value
$ws.Range("Q5").Value = 90.6
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)  # xlPasteFormats

# Match the author's resulting selection state
$ws.Range("P12").Select()
